# Israel Premier League workbook update (11-04-2024)
#  1. Row 16 and Row 17 had their match data (all columns except the
#     leading "id" column A) swapped.
#  2. Row 44 and Row 45 had their match data (all columns except the
#     leading "id" column A) swapped.
#  3. Four new fixtures were appended as rows 194-197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    # Swap every column from B (2) through AC (29); column A (the row's
    # "id" index) stays put on both rows.
    for ($c = 2; $c -le 29; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-RowData 16 17
Swap-RowData 44 45

function Add-Fixture($row, $id, $matchId, $dateSerial, $homeTeam, $awayTeam,
    $oddH_op, $oddD_op, $oddA_op, $oddH, $oddD, $oddA,
    $ah, $oddAHH, $oddAHA, $ahOU, $oddAHOver, $oddAHUnder) {

    # Copy formatting from the row immediately above (row 193 is a fully
    # populated, already-played fixture) so the new rows pick up the same
    # styles (bold/bordered id cell, date number format, ...).
    $ws.Range("A193").Copy()
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("E193").Copy()
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $id          # A id
    $ws.Cells.Item($row, 2).Value = $matchId     # B id (source match id)
    $ws.Cells.Item($row, 3).Value = "Israel Premier League"   # C Div
    $ws.Cells.Item($row, 4).Value = "Israel Premier League"   # D Div Original Name
    $ws.Cells.Item($row, 5).Value = $dateSerial  # E Date
    $ws.Cells.Item($row, 6).Value = $homeTeam    # F HomeTeam
    $ws.Cells.Item($row, 7).Value = $awayTeam    # G AwayTeam
    # H FTHG, I FTAG, J FTR intentionally left blank - not played yet.
    $ws.Cells.Item($row, 11).Value = $oddH_op    # K
    $ws.Cells.Item($row, 12).Value = $oddD_op    # L
    $ws.Cells.Item($row, 13).Value = $oddA_op    # M
    $ws.Cells.Item($row, 14).Value = $oddH       # N
    $ws.Cells.Item($row, 15).Value = $oddD       # O
    $ws.Cells.Item($row, 16).Value = $oddA       # P
    $ws.Cells.Item($row, 17).Value = $ah         # Q
    $ws.Cells.Item($row, 18).Value = $oddAHH     # R
    $ws.Cells.Item($row, 19).Value = $oddAHA     # S
    $ws.Cells.Item($row, 20).Value = $ahOU       # T
    $ws.Cells.Item($row, 21).Value = $oddAHOver  # U
    $ws.Cells.Item($row, 22).Value = $oddAHUnder # V
    $ws.Cells.Item($row, 23).Value = 0           # W PLH
    $ws.Cells.Item($row, 24).Value = 0           # X PLD
    $ws.Cells.Item($row, 25).Value = 0           # Y PLA
    $ws.Cells.Item($row, 26).Value = 0           # Z PL_Ahh
    $ws.Cells.Item($row, 27).Value = 0           # AA PL_Aha
    # AB PL_AhOver, AC PL_AhUnder intentionally left blank.
}

Add-Fixture 194 192 8015665 45395.48958333334 "Maccabi Haifa" "Hapoel Bnei Sakhnin" `
    1.4 4.5 7 1.363 4.75 7.5 -1.25 1.825 2.025 2.75 1.875 1.975

Add-Fixture 195 193 8016164 45395.58333333334 "Hapoel TelAviv" "Hapoel Petah Tikva" `
    2.3 3.2 3 2.3 3.2 3 -0.25 2.1 1.775 2.25 2.025 1.825

Add-Fixture 196 194 8015666 45396.57291666666 "Hapoel Haifa" "Maccabi Bnei Raina" `
    2.15 3.2 3.1 2.15 3.2 3.1 -0.25 1.925 1.925 2.25 1.975 1.875

Add-Fixture 197 195 8016159 45397.60416666666 "Hapoel Beer Sheva" "Maccabi Tel Aviv" `
    3.5 3.25 2 3.5 3.25 2 0.5 1.8 2.05 2.5 2.025 1.825

Write-Output "edit applied"
